# Apply actitud.xlsx edits:
#  - Astronauta: add column F ("S5") values for several students, recalculates L (Resultados)
#  - Senador: add column F ("S5") values for several students, recalculates L (Resultados)
#  - View state: Astronauta becomes the active/selected tab (cell F9 selected),
#    Senador scrolls to show row 2 at top with A26 selected, Mago is no longer
#    the active tab (selection stays at E11).

$wb = $excel.ActiveWorkbook

# --- Astronauta ---------------------------------------------------------
$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsAstronauta.Range("F4").Value = 1
$wsAstronauta.Range("F5").Value = 1
$wsAstronauta.Range("F8").Value = 1
$wsAstronauta.Range("F10").Value = 1
$wsAstronauta.Range("F15").Value = 1
$wsAstronauta.Range("F17").Value = 1

# --- Senador -------------------------------------------------------------
$wsSenador = $wb.Worksheets.Item("Senador")
$wsSenador.Range("F6").Value = 0
$wsSenador.Range("F9").Value = 0
$wsSenador.Range("F12").Value = 0
$wsSenador.Range("F18").Value = 0

# --- View / selection state ----------------------------------------------
# Senador: scroll so row 2 is at top, select A26 (Mago/Senador stay inactive
# at the end, so do this before activating Astronauta as the final tab).
$wsSenador.Activate()
$excel.ActiveWindow.ScrollRow = 2
$wsSenador.Range("A26").Select()

# Mago: keep its own selection as-is (E11); just make sure it is not the
# final active tab (handled by activating Astronauta last).
$wsMago = $wb.Worksheets.Item("Mago")
$wsMago.Activate()
$wsMago.Range("E11").Select()

# Astronauta: final active tab, with F9 selected.
$wsAstronauta.Activate()
$wsAstronauta.Range("F9").Select()
